$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update predicted-probability columns J:U (rows 2-61) with refreshed model output
$arr = New-Object "object[,]" 60,12
$arr[0,0] = 0.78233415
$arr[0,1] = 0.20511045
$arr[0,2] = 0.01255532
$arr[0,3] = 0.31529808
$arr[0,4] = 0.5278524
$arr[0,5] = 0.15684949
$arr[0,6] = 0.52108383
$arr[0,7] = 0.4559837
$arr[0,8] = 0.0229325
$arr[0,9] = 0.5807765
$arr[0,10] = 0.22852692
$arr[0,11] = 0.19069654
$arr[1,0] = 0.18075332
$arr[1,1] = 0.7986664
$arr[1,2] = 0.02058036
$arr[1,3] = 0.65408486
$arr[1,4] = 0.16951942
$arr[1,5] = 0.17639571
$arr[1,6] = 0.5706344
$arr[1,7] = 0.08771849
$arr[1,8] = 0.3416471
$arr[1,9] = 0.44118333
$arr[1,10] = 0.35157853
$arr[1,11] = 0.20723815
$arr[2,0] = 0.38078097
$arr[2,1] = 0.60347265
$arr[2,2] = 0.01574643
$arr[2,3] = 0.2924593
$arr[2,4] = 0.48588544
$arr[2,5] = 0.22165531
$arr[2,6] = 0.22798598
$arr[2,7] = 0.74200976
$arr[2,8] = 0.03000422
$arr[2,9] = 0.32608503
$arr[2,10] = 0.37868685
$arr[2,11] = 0.29522812
$arr[3,0] = 0.09594364
$arr[3,1] = 0.8728626
$arr[3,2] = 0.03119381
$arr[3,3] = 0.22517376
$arr[3,4] = 0.58907133
$arr[3,5] = 0.18575493
$arr[3,6] = 0.26230663
$arr[3,7] = 0.7064157
$arr[3,8] = 0.03127763
$arr[3,9] = 0.16058183
$arr[3,10] = 0.4987595
$arr[3,11] = 0.34065866
$arr[4,0] = 0.1045831
$arr[4,1] = 0.8553313
$arr[4,2] = 0.04008559
$arr[4,3] = 0.54612947
$arr[4,4] = 0.21363984
$arr[4,5] = 0.24023068
$arr[4,6] = 0.4155069
$arr[4,7] = 0.24809934
$arr[4,8] = 0.33639377
$arr[4,9] = 0.2623141
$arr[4,10] = 0.5814155
$arr[4,11] = 0.15627044
$arr[5,0] = 0.10935567
$arr[5,1] = 0.8194507
$arr[5,2] = 0.07119361
$arr[5,3] = 0.5091798
$arr[5,4] = 0.2530529
$arr[5,5] = 0.23776734
$arr[5,6] = 0.38257557
$arr[5,7] = 0.5793567
$arr[5,8] = 0.03806776
$arr[5,9] = 0.39482948
$arr[5,10] = 0.4015167
$arr[5,11] = 0.20365381
$arr[6,0] = 0.08724256
$arr[6,1] = 0.8371845
$arr[6,2] = 0.075573
$arr[6,3] = 0.22517376
$arr[6,4] = 0.58907133
$arr[6,5] = 0.18575493
$arr[6,6] = 0.2925673
$arr[6,7] = 0.68565595
$arr[6,8] = 0.02177675
$arr[6,9] = 0.23747084
$arr[6,10] = 0.56529033
$arr[6,11] = 0.19723883
$arr[7,0] = 0.13899976
$arr[7,1] = 0.8257836
$arr[7,2] = 0.03521656
$arr[7,3] = 0.6267757
$arr[7,4] = 0.18688317
$arr[7,5] = 0.18634117
$arr[7,6] = 0.2509495
$arr[7,7] = 0.51163507
$arr[7,8] = 0.23741543
$arr[7,9] = 0.4637538
$arr[7,10] = 0.35268748
$arr[7,11] = 0.18355876
$arr[8,0] = 0.73936427
$arr[8,1] = 0.20834522
$arr[8,2] = 0.05229051
$arr[8,3] = 0.3093528
$arr[8,4] = 0.4643321
$arr[8,5] = 0.22631513
$arr[8,6] = 0.46534243
$arr[8,7] = 0.5045809
$arr[8,8] = 0.03007668
$arr[8,9] = 0.4300297
$arr[8,10] = 0.3541671
$arr[8,11] = 0.21580324
$arr[9,0] = 0.05340498
$arr[9,1] = 0.6342567
$arr[9,2] = 0.31233832
$arr[9,3] = 0.22887723
$arr[9,4] = 0.5420357
$arr[9,5] = 0.22908705
$arr[9,6] = 0.3867906
$arr[9,7] = 0.5520421
$arr[9,8] = 0.06116728
$arr[9,9] = 0.21831052
$arr[9,10] = 0.580251
$arr[9,11] = 0.20143852
$arr[10,0] = 0.3621747
$arr[10,1] = 0.61812794
$arr[10,2] = 0.01969736
$arr[10,3] = 0.22896089
$arr[10,4] = 0.5821577
$arr[10,5] = 0.18888141
$arr[10,6] = 0.19232507
$arr[10,7] = 0.7781415
$arr[10,8] = 0.02953341
$arr[10,9] = 0.13979134
$arr[10,10] = 0.4246899
$arr[10,11] = 0.4355188
$arr[11,0] = 0.14209205
$arr[11,1] = 0.84466505
$arr[11,2] = 0.01324307
$arr[11,3] = 0.21030562
$arr[11,4] = 0.5701194
$arr[11,5] = 0.21957502
$arr[11,6] = 0.5589625
$arr[11,7] = 0.41401523
$arr[11,8] = 0.02702224
$arr[11,9] = 0.21885443
$arr[11,10] = 0.5804169
$arr[11,11] = 0.20072864
$arr[12,0] = 0.47218147
$arr[12,1] = 0.48456034
$arr[12,2] = 0.04325826
$arr[12,3] = 0.5554354
$arr[12,4] = 0.21510929
$arr[12,5] = 0.22945525
$arr[12,6] = 0.35211992
$arr[12,7] = 0.6124389
$arr[12,8] = 0.03544115
$arr[12,9] = 0.5735364
$arr[12,10] = 0.25054738
$arr[12,11] = 0.1759162
$arr[13,0] = 0.01777358
$arr[13,1] = 0.36885107
$arr[13,2] = 0.61337537
$arr[13,3] = 0.3611072
$arr[13,4] = 0.3742829
$arr[13,5] = 0.2646099
$arr[13,6] = 0.3743962
$arr[13,7] = 0.39082798
$arr[13,8] = 0.23477581
$arr[13,9] = 0.21177322
$arr[13,10] = 0.48953098
$arr[13,11] = 0.29869583
$arr[14,0] = 0.55213344
$arr[14,1] = 0.2797109
$arr[14,2] = 0.16815573
$arr[14,3] = 0.41474745
$arr[14,4] = 0.2673323
$arr[14,5] = 0.31792027
$arr[14,6] = 0.2131971
$arr[14,7] = 0.6561246
$arr[14,8] = 0.13067827
$arr[14,9] = 0.52597195
$arr[14,10] = 0.23405682
$arr[14,11] = 0.23997125
$arr[15,0] = 0.11041693
$arr[15,1] = 0.87948173
$arr[15,2] = 0.01010134
$arr[15,3] = 0.22517376
$arr[15,4] = 0.58907133
$arr[15,5] = 0.18575493
$arr[15,6] = 0.68055975
$arr[15,7] = 0.2850601
$arr[15,8] = 0.03438013
$arr[15,9] = 0.43984276
$arr[15,10] = 0.38366708
$arr[15,11] = 0.17649013
$arr[16,0] = 0.12549646
$arr[16,1] = 0.8354964
$arr[16,2] = 0.03900701
$arr[16,3] = 0.27403829
$arr[16,4] = 0.4998965
$arr[16,5] = 0.22606523
$arr[16,6] = 0.10710905
$arr[16,7] = 0.86917627
$arr[16,8] = 0.02371471
$arr[16,9] = 0.28393844
$arr[16,10] = 0.4952324
$arr[16,11] = 0.22082913
$arr[17,0] = 0.6986749
$arr[17,1] = 0.25418133
$arr[17,2] = 0.04714379
$arr[17,3] = 0.5878008
$arr[17,4] = 0.21540059
$arr[17,5] = 0.19679864
$arr[17,6] = 0.8351672
$arr[17,7] = 0.15179847
$arr[17,8] = 0.01303429
$arr[17,9] = 0.587035
$arr[17,10] = 0.22649983
$arr[17,11] = 0.18646516
$arr[18,0] = 0.24233446
$arr[18,1] = 0.6030499
$arr[18,2] = 0.1546158
$arr[18,3] = 0.3873328
$arr[18,4] = 0.44024625
$arr[18,5] = 0.17242095
$arr[18,6] = 0.84069496
$arr[18,7] = 0.033853
$arr[18,8] = 0.12545201
$arr[18,9] = 0.5714529
$arr[18,10] = 0.23696268
$arr[18,11] = 0.19158444
$arr[19,0] = 0.65966445
$arr[19,1] = 0.2713271
$arr[19,2] = 0.06900843
$arr[19,3] = 0.4273544
$arr[19,4] = 0.3064588
$arr[19,5] = 0.26618677
$arr[19,6] = 0.74247324
$arr[19,7] = 0.22158456
$arr[19,8] = 0.03594225
$arr[19,9] = 0.47960028
$arr[19,10] = 0.33378404
$arr[19,11] = 0.18661566
$arr[20,0] = 0.03846259
$arr[20,1] = 0.8997402
$arr[20,2] = 0.06179722
$arr[20,3] = 0.57133657
$arr[20,4] = 0.28600737
$arr[20,5] = 0.14265604
$arr[20,6] = 0.18387789
$arr[20,7] = 0.7231051
$arr[20,8] = 0.09301709
$arr[20,9] = 0.20098317
$arr[20,10] = 0.596818
$arr[20,11] = 0.20219892
$arr[21,0] = 0.50567704
$arr[21,1] = 0.48791695
$arr[21,2] = 0.00640599
$arr[21,3] = 0.64132017
$arr[21,4] = 0.1901995
$arr[21,5] = 0.16848038
$arr[21,6] = 0.74665827
$arr[21,7] = 0.18384798
$arr[21,8] = 0.06949377
$arr[21,9] = 0.50291735
$arr[21,10] = 0.291104
$arr[21,11] = 0.20597866
$arr[22,0] = 0.19513567
$arr[22,1] = 0.72549325
$arr[22,2] = 0.07937109
$arr[22,3] = 0.3450672
$arr[22,4] = 0.43559083
$arr[22,5] = 0.21934193
$arr[22,6] = 0.36279303
$arr[22,7] = 0.3675772
$arr[22,8] = 0.26962978
$arr[22,9] = 0.392653
$arr[22,10] = 0.3742651
$arr[22,11] = 0.23308189
$arr[23,0] = 0.11909114
$arr[23,1] = 0.8687427
$arr[23,2] = 0.01216615
$arr[23,3] = 0.1641108
$arr[23,4] = 0.50934696
$arr[23,5] = 0.32654223
$arr[23,6] = 0.17051364
$arr[23,7] = 0.4893998
$arr[23,8] = 0.34008655
$arr[23,9] = 0.23258983
$arr[23,10] = 0.24842677
$arr[23,11] = 0.5189835
$arr[24,0] = 0.30455458
$arr[24,1] = 0.6088488
$arr[24,2] = 0.08659661
$arr[24,3] = 0.35861188
$arr[24,4] = 0.41371983
$arr[24,5] = 0.22766834
$arr[24,6] = 0.11157916
$arr[24,7] = 0.66456556
$arr[24,8] = 0.22385527
$arr[24,9] = 0.4020642
$arr[24,10] = 0.39845574
$arr[24,11] = 0.19948004
$arr[25,0] = 0.16643424
$arr[25,1] = 0.8268189
$arr[25,2] = 0.0067468
$arr[25,3] = 0.29928565
$arr[25,4] = 0.2786294
$arr[25,5] = 0.422085
$arr[25,6] = 0.8533324
$arr[25,7] = 0.08835965
$arr[25,8] = 0.05830795
$arr[25,9] = 0.49178234
$arr[25,10] = 0.33512753
$arr[25,11] = 0.17309016
$arr[26,0] = 0.0854371
$arr[26,1] = 0.60524815
$arr[26,2] = 0.30931482
$arr[26,3] = 0.5151216
$arr[26,4] = 0.25130117
$arr[26,5] = 0.23357724
$arr[26,6] = 0.24798535
$arr[26,7] = 0.6242253
$arr[26,8] = 0.12778936
$arr[26,9] = 0.522357
$arr[26,10] = 0.27823472
$arr[26,11] = 0.19940834
$arr[27,0] = 0.12353331
$arr[27,1] = 0.5679025
$arr[27,2] = 0.30856416
$arr[27,3] = 0.35348335
$arr[27,4] = 0.3787492
$arr[27,5] = 0.26776746
$arr[27,6] = 0.5289345
$arr[27,7] = 0.44908234
$arr[27,8] = 0.0219832
$arr[27,9] = 0.29444182
$arr[27,10] = 0.44910902
$arr[27,11] = 0.25644913
$arr[28,0] = 0.21167032
$arr[28,1] = 0.758362
$arr[28,2] = 0.0299677
$arr[28,3] = 0.39214933
$arr[28,4] = 0.3584595
$arr[28,5] = 0.24939118
$arr[28,6] = 0.29644194
$arr[28,7] = 0.6751565
$arr[28,8] = 0.02840157
$arr[28,9] = 0.6354073
$arr[28,10] = 0.20118864
$arr[28,11] = 0.16340408
$arr[29,0] = 0.3490355
$arr[29,1] = 0.64091486
$arr[29,2] = 0.01004969
$arr[29,3] = 0.19402437
$arr[29,4] = 0.5924538
$arr[29,5] = 0.21352185
$arr[29,6] = 0.425008
$arr[29,7] = 0.09320848
$arr[29,8] = 0.4817835
$arr[29,9] = 0.5901689
$arr[29,10] = 0.25206605
$arr[29,11] = 0.15776503
$arr[30,0] = 0.73696226
$arr[30,1] = 0.20762178
$arr[30,2] = 0.055416
$arr[30,3] = 0.6208093
$arr[30,4] = 0.178932
$arr[30,5] = 0.20025867
$arr[30,6] = 0.9327992
$arr[30,7] = 0.05786263
$arr[30,8] = 0.00933817
$arr[30,9] = 0.6317823
$arr[30,10] = 0.2028209
$arr[30,11] = 0.16539681
$arr[31,0] = 0.12473991
$arr[31,1] = 0.850543
$arr[31,2] = 0.02471699
$arr[31,3] = 0.64347696
$arr[31,4] = 0.1815266
$arr[31,5] = 0.17499644
$arr[31,6] = 0.5262218
$arr[31,7] = 0.4276661
$arr[31,8] = 0.04611213
$arr[31,9] = 0.39057767
$arr[31,10] = 0.440992
$arr[31,11] = 0.16843033
$arr[32,0] = 0.4103611
$arr[32,1] = 0.49737862
$arr[32,2] = 0.09226028
$arr[32,3] = 0.43062958
$arr[32,4] = 0.3573999
$arr[32,5] = 0.2119705
$arr[32,6] = 0.20505385
$arr[32,7] = 0.7610595
$arr[32,8] = 0.0338866
$arr[32,9] = 0.44685352
$arr[32,10] = 0.21720336
$arr[32,11] = 0.33594307
$arr[33,0] = 0.22173472
$arr[33,1] = 0.7632047
$arr[33,2] = 0.01506061
$arr[33,3] = 0.23814596
$arr[33,4] = 0.5659863
$arr[33,5] = 0.19586775
$arr[33,6] = 0.50611264
$arr[33,7] = 0.44264364
$arr[33,8] = 0.05124369
$arr[33,9] = 0.30765206
$arr[33,10] = 0.509218
$arr[33,11] = 0.18312998
$arr[34,0] = 0.58734137
$arr[34,1] = 0.37736794
$arr[34,2] = 0.03529063
$arr[34,3] = 0.58672065
$arr[34,4] = 0.21382213
$arr[34,5] = 0.19945726
$arr[34,6] = 0.7173649
$arr[34,7] = 0.21917927
$arr[34,8] = 0.06345579
$arr[34,9] = 0.4965223
$arr[34,10] = 0.34737575
$arr[34,11] = 0.15610191
$arr[35,0] = 0.09024981
$arr[35,1] = 0.90366215
$arr[35,2] = 0.00608824
$arr[35,3] = 0.5607542
$arr[35,4] = 0.3004313
$arr[35,5] = 0.13881452
$arr[35,6] = 0.9421651
$arr[35,7] = 0.04636282
$arr[35,8] = 0.01147209
$arr[35,9] = 0.571793
$arr[35,10] = 0.28241724
$arr[35,11] = 0.1457898
$arr[36,0] = 0.46739677
$arr[36,1] = 0.47537032
$arr[36,2] = 0.05723293
$arr[36,3] = 0.4090801
$arr[36,4] = 0.26922634
$arr[36,5] = 0.32169357
$arr[36,6] = 0.37088436
$arr[36,7] = 0.19510153
$arr[36,8] = 0.43401408
$arr[36,9] = 0.457657
$arr[36,10] = 0.3532945
$arr[36,11] = 0.1890485
$arr[37,0] = 0.7790382
$arr[37,1] = 0.15656023
$arr[37,2] = 0.06440157
$arr[37,3] = 0.63841903
$arr[37,4] = 0.17080575
$arr[37,5] = 0.19077516
$arr[37,6] = 0.84987897
$arr[37,7] = 0.09881976
$arr[37,8] = 0.05130132
$arr[37,9] = 0.63232553
$arr[37,10] = 0.19240321
$arr[37,11] = 0.17527121
$arr[38,0] = 0.359148
$arr[38,1] = 0.5845576
$arr[38,2] = 0.05629439
$arr[38,3] = 0.19214351
$arr[38,4] = 0.61453915
$arr[38,5] = 0.19331735
$arr[38,6] = 0.44586864
$arr[38,7] = 0.5185939
$arr[38,8] = 0.03553742
$arr[38,9] = 0.29229927
$arr[38,10] = 0.45128745
$arr[38,11] = 0.25641328
$arr[39,0] = 0.04453562
$arr[39,1] = 0.94974715
$arr[39,2] = 0.00571728
$arr[39,3] = 0.17227101
$arr[39,4] = 0.64788467
$arr[39,5] = 0.17984432
$arr[39,6] = 0.88462704
$arr[39,7] = 0.06544604
$arr[39,8] = 0.04992694
$arr[39,9] = 0.34678522
$arr[39,10] = 0.44272405
$arr[39,11] = 0.21049073
$arr[40,0] = 0.3239441
$arr[40,1] = 0.6704308
$arr[40,2] = 0.00562515
$arr[40,3] = 0.57513714
$arr[40,4] = 0.2366909
$arr[40,5] = 0.188172
$arr[40,6] = 0.57850754
$arr[40,7] = 0.40917075
$arr[40,8] = 0.0123217
$arr[40,9] = 0.6074341
$arr[40,10] = 0.2198658
$arr[40,11] = 0.17270014
$arr[41,0] = 0.01234764
$arr[41,1] = 0.97458935
$arr[41,2] = 0.01306298
$arr[41,3] = 0.5644657
$arr[41,4] = 0.17074728
$arr[41,5] = 0.26478702
$arr[41,6] = 0.1839202
$arr[41,7] = 0.6074712
$arr[41,8] = 0.20860857
$arr[41,9] = 0.2334138
$arr[41,10] = 0.49437046
$arr[41,11] = 0.27221575
$arr[42,0] = 0.76620567
$arr[42,1] = 0.22507893
$arr[42,2] = 0.00871531
$arr[42,3] = 0.31325087
$arr[42,4] = 0.30848697
$arr[42,5] = 0.3782622
$arr[42,6] = 0.304784
$arr[42,7] = 0.56788576
$arr[42,8] = 0.12733026
$arr[42,9] = 0.24926138
$arr[42,10] = 0.46809393
$arr[42,11] = 0.2826447
$arr[43,0] = 0.44310105
$arr[43,1] = 0.5502323
$arr[43,2] = 0.00666667
$arr[43,3] = 0.5934349
$arr[43,4] = 0.18264559
$arr[43,5] = 0.22391956
$arr[43,6] = 0.9095509
$arr[43,7] = 0.05617945
$arr[43,8] = 0.03426965
$arr[43,9] = 0.59881294
$arr[43,10] = 0.24063183
$arr[43,11] = 0.16055524
$arr[44,0] = 0.200317
$arr[44,1] = 0.7828024
$arr[44,2] = 0.01688055
$arr[44,3] = 0.6427025
$arr[44,4] = 0.18698886
$arr[44,5] = 0.17030863
$arr[44,6] = 0.8157881
$arr[44,7] = 0.13121767
$arr[44,8] = 0.05299423
$arr[44,9] = 0.57345223
$arr[44,10] = 0.2785115
$arr[44,11] = 0.14803627
$arr[45,0] = 0.39605165
$arr[45,1] = 0.5057923
$arr[45,2] = 0.09815606
$arr[45,3] = 0.45709538
$arr[45,4] = 0.36209399
$arr[45,5] = 0.18081065
$arr[45,6] = 0.21960452
$arr[45,7] = 0.7469027
$arr[45,8] = 0.03349279
$arr[45,9] = 0.34086186
$arr[45,10] = 0.44206935
$arr[45,11] = 0.21706876
$arr[46,0] = 0.68862426
$arr[46,1] = 0.21345183
$arr[46,2] = 0.0979239
$arr[46,3] = 0.652636
$arr[46,4] = 0.17591082
$arr[46,5] = 0.17145315
$arr[46,6] = 0.9683416
$arr[46,7] = 0.02305705
$arr[46,8] = 0.00860133
$arr[46,9] = 0.64035165
$arr[46,10] = 0.19965121
$arr[46,11] = 0.15999717
$arr[47,0] = 0.6010738
$arr[47,1] = 0.30029324
$arr[47,2] = 0.09863294
$arr[47,3] = 0.6508088
$arr[47,4] = 0.17488426
$arr[47,5] = 0.17430691
$arr[47,6] = 0.93630147
$arr[47,7] = 0.05333555
$arr[47,8] = 0.01036297
$arr[47,9] = 0.5389509
$arr[47,10] = 0.25359318
$arr[47,11] = 0.20745589
$arr[48,0] = 0.6286804
$arr[48,1] = 0.36106598
$arr[48,2] = 0.01025364
$arr[48,3] = 0.6508088
$arr[48,4] = 0.17488426
$arr[48,5] = 0.17430691
$arr[48,6] = 0.563457
$arr[48,7] = 0.41129112
$arr[48,8] = 0.02525187
$arr[48,9] = 0.5976043
$arr[48,10] = 0.250569
$arr[48,11] = 0.15182671
$arr[49,0] = 0.5097749
$arr[49,1] = 0.47763667
$arr[49,2] = 0.01258841
$arr[49,3] = 0.5613158
$arr[49,4] = 0.2716416
$arr[49,5] = 0.16704261
$arr[49,6] = 0.72395074
$arr[49,7] = 0.23810215
$arr[49,8] = 0.03794715
$arr[49,9] = 0.6324192
$arr[49,10] = 0.20803878
$arr[49,11] = 0.15954198
$arr[50,0] = 0.06529194
$arr[50,1] = 0.92276686
$arr[50,2] = 0.01194129
$arr[50,3] = 0.23814596
$arr[50,4] = 0.5659863
$arr[50,5] = 0.19586775
$arr[50,6] = 0.21870007
$arr[50,7] = 0.7658943
$arr[50,8] = 0.0154056
$arr[50,9] = 0.48498422
$arr[50,10] = 0.3355346
$arr[50,11] = 0.17948115
$arr[51,0] = 0.683288
$arr[51,1] = 0.27160034
$arr[51,2] = 0.04511166
$arr[51,3] = 0.6508088
$arr[51,4] = 0.17488426
$arr[51,5] = 0.17430691
$arr[51,6] = 0.44216996
$arr[51,7] = 0.5438422
$arr[51,8] = 0.01398788
$arr[51,9] = 0.57958305
$arr[51,10] = 0.27405536
$arr[51,11] = 0.14636156
$arr[52,0] = 0.6375044
$arr[52,1] = 0.33775854
$arr[52,2] = 0.02473706
$arr[52,3] = 0.44337472
$arr[52,4] = 0.37268606
$arr[52,5] = 0.18393925
$arr[52,6] = 0.7675657
$arr[52,7] = 0.21500853
$arr[52,8] = 0.0174257
$arr[52,9] = 0.58743054
$arr[52,10] = 0.2644748
$arr[52,11] = 0.14809465
$arr[53,0] = 0.0383087
$arr[53,1] = 0.90962756
$arr[53,2] = 0.05206389
$arr[53,3] = 0.23814596
$arr[53,4] = 0.5659863
$arr[53,5] = 0.19586775
$arr[53,6] = 0.5786059
$arr[53,7] = 0.38503477
$arr[53,8] = 0.03635934
$arr[53,9] = 0.22674003
$arr[53,10] = 0.5708306
$arr[53,11] = 0.20242938
$arr[54,0] = 0.687098
$arr[54,1] = 0.21986641
$arr[54,2] = 0.09303556
$arr[54,3] = 0.43296948
$arr[54,4] = 0.25610837
$arr[54,5] = 0.31092212
$arr[54,6] = 0.37371668
$arr[54,7] = 0.2961685
$arr[54,8] = 0.3301148
$arr[54,9] = 0.42483515
$arr[54,10] = 0.3566059
$arr[54,11] = 0.21855897
$arr[55,0] = 0.09903228
$arr[55,1] = 0.860985
$arr[55,2] = 0.03998286
$arr[55,3] = 0.240357
$arr[55,4] = 0.56093854
$arr[55,5] = 0.19870444
$arr[55,6] = 0.08647666
$arr[55,7] = 0.8912964
$arr[55,8] = 0.02222694
$arr[55,9] = 0.22922014
$arr[55,10] = 0.46041176
$arr[55,11] = 0.3103681
$arr[56,0] = 0.5553298
$arr[56,1] = 0.43692684
$arr[56,2] = 0.00774332
$arr[56,3] = 0.6508088
$arr[56,4] = 0.17488426
$arr[56,5] = 0.17430691
$arr[56,6] = 0.76533365
$arr[56,7] = 0.22077835
$arr[56,8] = 0.013888
$arr[56,9] = 0.64643705
$arr[56,10] = 0.19305614
$arr[56,11] = 0.16050678
$arr[57,0] = 0.12933844
$arr[57,1] = 0.85998183
$arr[57,2] = 0.01067977
$arr[57,3] = 0.4321422
$arr[57,4] = 0.37742636
$arr[57,5] = 0.19043149
$arr[57,6] = 0.19466281
$arr[57,7] = 0.70861113
$arr[57,8] = 0.09672609
$arr[57,9] = 0.25741792
$arr[57,10] = 0.50232625
$arr[57,11] = 0.24025589
$arr[58,0] = 0.3084658
$arr[58,1] = 0.6773461
$arr[58,2] = 0.01418821
$arr[58,3] = 0.56053805
$arr[58,4] = 0.24833256
$arr[58,5] = 0.19112942
$arr[58,6] = 0.87928736
$arr[58,7] = 0.08407827
$arr[58,8] = 0.03663441
$arr[58,9] = 0.52673197
$arr[58,10] = 0.2922024
$arr[58,11] = 0.1810657
$arr[59,0] = 0.12728818
$arr[59,1] = 0.8481103
$arr[59,2] = 0.02460149
$arr[59,3] = 0.21762916
$arr[59,4] = 0.58324546
$arr[59,5] = 0.19912538
$arr[59,6] = 0.33769798
$arr[59,7] = 0.6034233
$arr[59,8] = 0.05887873
$arr[59,9] = 0.42277628
$arr[59,10] = 0.40801203
$arr[59,11] = 0.16921176
$ws.Range("J2:U61").Value2 = $arr

# Move the active selection to U3 (matches the author's final cursor position)
$ws.Range("U3").Select()
